$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 (Objetivos:): replace the long PT objective text with the docente
# identification string (this text later reappears lower in the sheet too).
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C10").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"

# ---------------------------------------------------------------------------
# Row 13 used to be an orphan docente row (B/C only). It becomes the
# "Programa resumido:" row, now holding "Semestral".
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------------
# Row 14 becomes "Short syllabus:" with the (unchanged) English short
# syllabus description text.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Descrição do programa resumido em inglês.Real numbers, real functions, limits and derivatives of real functions.  Applications of the derivative and Taylor’s Formula."
$ws.Range("C14").Value = "Descrição do programa resumido em inglês.Real numbers, real functions, limits and derivatives of real functions.  Applications of the derivative and Taylor’s Formula."

# ---------------------------------------------------------------------------
# Row 15 becomes "Programa:" with B/C holding the literal text "01/01/2018"
# (a stray re-use of the same text as the "Ativação:" row). Entering that
# text directly would be auto-parsed into a date serial by Excel, so instead
# copy the already-correct text+style from B8:C8 (the real "Ativação:" row,
# which already stores that exact string as text) onto B15:C15.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8:C8").Copy()
$ws.Range("B15:C15").PasteSpecial(-4104)
$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------------------
# Row 16 becomes "Syllabus:" with the (unchanged) English syllabus text.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "•Real Numbers and Real Functions: trigonometric, exponential and logarithmic functions. Composite and inverse functions.•Limits: Definition, algebraic properties and squeeze theorem. Infinite limits and Limits to infinite.•Continuity: Weierstrass theorem and intermediate value theorem.•Derivative of real functions: Definition, geometrical and physics interpretations, derivative rules, chain rule, derivative of inverse and implicit functions, l’hopital rule, mean value theorem and consequences, Taylor’s Formula,  Maximum and Minimum Problems"
$ws.Range("C16").Value = "•Real Numbers and Real Functions: trigonometric, exponential and logarithmic functions. Composite and inverse functions.•Limits: Definition, algebraic properties and squeeze theorem. Infinite limits and Limits to infinite.•Continuity: Weierstrass theorem and intermediate value theorem.•Derivative of real functions: Definition, geometrical and physics interpretations, derivative rules, chain rule, derivative of inverse and implicit functions, l’hopital rule, mean value theorem and consequences, Taylor’s Formula,  Maximum and Minimum Problems"

# ---------------------------------------------------------------------------
# Row 17 becomes "Avaliação:" — a label-only row (no B/C values), matching
# the other label-only rows such as "Docentes responsáveis:".
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).RowHeight = 15

# ---------------------------------------------------------------------------
# Row 18 becomes "Método:" and now carries the docente string again in B/C.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Range("C18").Value = "5840692 - Diovana Aparecida dos Santos Napoleão"
$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------------------
# Rows 19-21 just shift their labels up by one slot; the B/C contents (the
# grading-criteria paragraphs) stay put.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------------
# The old row 22 (Bibliografia: / the STEWART... reference list) is removed
# entirely, shrinking the sheet to A1:C21.
# ---------------------------------------------------------------------------
$ws.Rows.Item(22).Delete()
